# Generate Report for Handback
#
# This localization-status workbook has three sheets:
#   Overview (sheet1), zh-cn (sheet2), de-de (sheet3)
#
# The handback run completed for both zh-cn and de-de:
#   - Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears
#     (Overview!B2:C3, zh-cn!C2:C3, de-de!C2:C3).
#   - Each language sheet grows two new columns describing the
#     handback itself: F "Latest Target File" and G "Latest Handback
#     File", both populated + hyperlinked for row 2 and row 3.
#   - Column H "Latest Handback DateTime" moves off the
#     "0001-01-01 00:00:00" placeholder and on to the real
#     handback timestamp for each language.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status column for both languages, both rows ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$zhAMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4c74abd674fc871348c97f2c05f5080a03bf776e/e2e/a.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bcc10210a64f24377052fa357d051284cfac8052/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhHandbackTime = "2016-03-19 22:26:30"

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# New "Latest Target File" (F) and "Latest Handback File" (G) columns,
# row 2 and row 3 - both rows point at the same handed-back file.
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhAMdUrl, "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhAMdUrl, "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

# "Latest Handback DateTime" (H) - was the 0001-01-01 placeholder.
$wsZhCn.Range("H2").Value = $zhHandbackTime
$wsZhCn.Range("H3").Value = $zhHandbackTime

# --- de-de sheet ---
$deAMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4c74abd674fc871348c97f2c05f5080a03bf776e/e2e/a.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b65050e0c9c727c3b6950313de6a82a833ced651/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deHandbackTime = "2016-03-19 22:26:36"

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deAMdUrl, "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deAMdUrl, "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$wsDeDe.Range("H2").Value = $deHandbackTime
$wsDeDe.Range("H3").Value = $deHandbackTime
